$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.689.17'
$ws.Range("E2").Value = '  -9.02%  '
$ws.Range("D3").Value = '2.421.03'
$ws.Range("E3").Value = '  -11.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '462.90'
$ws.Range("E5").Value = '  -8.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.81'
$ws.Range("E6").Value = '  -7.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.485'
$ws.Range("E8").Value = '  -9.32%  '
$ws.Range("D9").Value = '2.424.29'
$ws.Range("E9").Value = '  -11.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0949'
$ws.Range("E10").Value = '  -9.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.32'
$ws.Range("E11").Value = '  -12.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.314'
$ws.Range("E12").Value = '  -10.86%  '
$ws.Range("E13").Value = '  -4.22%  '
$ws.Range("D14").Value = '2.843.47'
$ws.Range("E14").Value = '  -11.94%  '
$ws.Range("D15").Value = '53.731.32'
$ws.Range("E15").Value = '  -9.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000133'
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.66'
$ws.Range("E17").Value = '  -10.04%  '
$ws.Range("D18").Value = '2.417.23'
$ws.Range("E18").Value = '  -12.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.18'
$ws.Range("E19").Value = '  -12.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '307.24'
$ws.Range("E20").Value = '  -11.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.43'
$ws.Range("E21").Value = '  -14.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.34'
$ws.Range("E24").Value = '  -15.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '55.97'
$ws.Range("E25").Value = '  -11.95%  '
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.386'
$ws.Range("E27").Value = '  -10.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.154'
$ws.Range("E28").Value = '  -11.64%  '
$ws.Range("D29").Value = '2.506.03'
$ws.Range("E29").Value = '  -12.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  -5.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").Value = '0.0₃0719'
$ws.Range("E32").Value = '  -14.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.03'
$ws.Range("E33").Value = '  -2.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.68'
$ws.Range("E34").Value = '  -8.46%  '
$ws.Range("E35").Value = '  -11.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.98'
$ws.Range("E36").Value = '  -8.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.52'
$ws.Range("E37").Value = '  -16.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").Value = '  -6.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.802'
$ws.Range("E39").Value = '  -16.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.995'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.79'
$ws.Range("E41").Value = '  -9.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.594'
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0525'
$ws.Range("E43").Value = '  -6.58%  '
$ws.Range("E44").Value = '  -8.67%  '
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("E46").Value = '  -12.25%  '
$ws.Range("D47").Value = '1.934.59'
$ws.Range("E47").Value = '  -11.61%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0870'
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0217'
$ws.Range("E49").Value = '  -4.55%  '
$ws.Range("E50").Value = '  -13.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.46'
$ws.Range("E51").Value = '  -14.31%  '
